$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The "DATA DICTIONARY" table occupies B192:H206 (header row 192 + 14 data
# rows). A new "Null" column is being inserted as column F (between the
# existing "Data Format" and "Field Size" columns), pushing the old F,G,H
# columns one place to the right, to become G,H,I.
#
# Strategy:
#  1. Shift the *formatting* (cell style) of F:H one column to the right,
#     processing right-to-left (H->I, G->H, F->G) so the in-progress shift
#     never reads a column that was already overwritten.
#  2. Give the new F column (rows 192:206) the same style as column C's
#     header/body cells (s=8 for the header row, s=9 for the data rows).
#  3. Write the final values straight into F,G,H,I for every row, row by
#     row, top to bottom, left to right - this also makes new shared
#     strings get interned in the same order Excel would have produced
#     them (Null, NOT NULL, NULL, Date).
# ---------------------------------------------------------------------------

# --- 1. shift existing formatting right: H->I, G->H, F->G ------------------
$ws.Range("H192:H206").Copy()
$ws.Range("I192:I206").PasteSpecial(-4122)

$ws.Range("G192:G206").Copy()
$ws.Range("H192:H206").PasteSpecial(-4122)

$ws.Range("F192:F206").Copy()
$ws.Range("G192:G206").PasteSpecial(-4122)

# --- 2. style the brand new column F ---------------------------------------
$ws.Range("C192").Copy()
$ws.Range("F192").PasteSpecial(-4122)

$ws.Range("C193").Copy()
$ws.Range("F193:F206").PasteSpecial(-4122)

# --- 3. write the final cell values -----------------------------------------

# Header row
$ws.Range("F192").Value = "Null"
$ws.Range("G192").Value = "Field Size"
$ws.Range("H192").Value = "Description"
$ws.Range("I192").Value = "Example"

# Customer.CustomerID
$ws.Range("F193").Value = "NOT NULL"
$ws.Range("G193").Value = 5
$ws.Range("H193").Value = "Unique ID for customer"
$ws.Range("I193").Value = 3

# Customer.CustomerFirstName
$ws.Range("F194").Value = "NULL"
$ws.Range("G194").Value = 20
$ws.Range("H194").Value = "First name for customer"
$ws.Range("I194").Value = "John"

# Customer.CustomerSurname
$ws.Range("F195").Value = "NULL"
$ws.Range("G195").Value = 20
$ws.Range("H195").Value = "Surname for customer"
$ws.Range("I195").Value = "Smith"

# Order.OrderID
$ws.Range("F196").Value = "NOT NULL"
$ws.Range("G196").Value = 5
$ws.Range("H196").Value = "Unique ID for order"
$ws.Range("I196").Value = 4

# Order.OrderDate (Data Type Datetime -> Date)
$ws.Range("D197").Value = "Date"
$ws.Range("F197").Value = "NOT NULL"
$ws.Range("G197").Value = 10
$ws.Range("H197").Value = "Date of order"
$ws.Range("I197").Value = 43575

# Order.OrderTotalPrice
$ws.Range("F198").Value = "NULL"
$ws.Range("G198").Value = 5
$ws.Range("H198").Value = "Total price of order"
$ws.Range("I198").Value = 13.99

# OrderDetails.ProductOrderQuantity
$ws.Range("F199").Value = "NOT NULL"
$ws.Range("G199").Value = 2
$ws.Range("H199").Value = "Amount of item ordered"
$ws.Range("I199").Value = 5

# OrderDetails.ProductOrderPrice
$ws.Range("F200").Value = "NULL"
$ws.Range("G200").Value = 5
$ws.Range("H200").Value = "Price of product * quantity"
$ws.Range("I200").Value = 5.99

# Product.ProductID
$ws.Range("F201").Value = "NOT NULL"
$ws.Range("G201").Value = 5
$ws.Range("H201").Value = "Unique ID for product"
$ws.Range("I201").Value = 6

# Product.ProductType
$ws.Range("F202").Value = "NULL"
$ws.Range("G202").Value = 5
$ws.Range("H202").Value = "Value for type of product"
$ws.Range("I202").Value = "Drink"

# Product.ProductPrice
$ws.Range("F203").Value = "NULL"
$ws.Range("G203").Value = 5
$ws.Range("H203").Value = "Price of individual product"
$ws.Range("I203").Value = 3.55

# Product.ProductName
$ws.Range("F204").Value = "NULL"
$ws.Range("G204").Value = 20
$ws.Range("H204").Value = "Name of product"
$ws.Range("I204").Value = "Croissant"

# Product.ProductCalories
$ws.Range("F205").Value = "NULL"
$ws.Range("G205").Value = 4
$ws.Range("H205").Value = "Calories of product"
$ws.Range("I205").Value = 400

# Product.ProductDetails
$ws.Range("F206").Value = "NULL"
$ws.Range("G206").Value = "max"
$ws.Range("H206").Value = "Details of product"
$ws.Range("I206").Value = "buttery, flaky pastry"

# ---------------------------------------------------------------------------
# Column widths: B:I now share one uniform width (30.71), same as the old
# B:G width group, J keeps its own width, column A/rest stay default.
# ---------------------------------------------------------------------------
$ws.Range("B1:I1").EntireColumn.ColumnWidth = 30.7109375

# ---------------------------------------------------------------------------
# View state: scrolled a bit further down, selection moved to G184.
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 177
$ws.Range("G184").Select()

Write-Output "done"
